# 自动更新Excel文件 - 2025-11-22 23:12:25
# For every data row (row 2..99) on the active sheet, column E holds the
# number of remaining days and column F holds the start date (yyyymmdd).
# On a daily refresh, remaining days (E) normally decreases by 1.
# Rows whose remaining days already reached 1 are treated as renewed:
# E is reset back up to 10 and F (start date) is advanced to the new
# "today" (20251123).
#
# Row 36 is intentionally left untouched because its F value is a
# malformed date (202510929) in the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$renewRows = @(50, 51, 52, 53, 54, 55, 56, 57)

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }

    if ($renewRows -contains $row) {
        $ws.Cells.Item($row, 5).Value2 = 10
        $ws.Cells.Item($row, 6).Value2 = 20251123
    } else {
        $current = $ws.Cells.Item($row, 5).Value2
        $ws.Cells.Item($row, 5).Value2 = $current - 1
    }
}
